$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting (style) of the last existing data row (29) down into the
# two new rows (30 and 31) so the new cells share the same style index as the
# rest of the data table.
$ws.Range($ws.Cells.Item(29, 1), $ws.Cells.Item(29, 18)).Copy()
$ws.Range($ws.Cells.Item(30, 1), $ws.Cells.Item(31, 18)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 30: Novokhopyorsk municipal district, female, year 2012
$ws.Cells.Item(30, 1).Value = 20627000
$ws.Cells.Item(30, 2).Value = "Новохопёрский муниципальный район"
$ws.Cells.Item(30, 3).Value = "female"
$ws.Cells.Item(30, 4).Value = 2012
$ws.Cells.Item(30, 5).Value = 0.05957
$ws.Cells.Item(30, 6).Value = 0.05304
$ws.Cells.Item(30, 7).Value = 0.0557
$ws.Cells.Item(30, 8).Value = 0.05145
$ws.Cells.Item(30, 9).Value = 0.06537
$ws.Cells.Item(30, 10).Value = 0.0779
$ws.Cells.Item(30, 11).Value = 0.0769
$ws.Cells.Item(30, 12).Value = 0.077
$ws.Cells.Item(30, 13).Value = 0.075
$ws.Cells.Item(30, 14).Value = 0.0861
$ws.Cells.Item(30, 15).Value = 0.10284
$ws.Cells.Item(30, 16).Value = 0.0891
$ws.Cells.Item(30, 17).Value = 0.0856
$ws.Cells.Item(30, 18).Value = 0.0444

# Row 31: Novokhopyorsk municipal district, male, year 2012
$ws.Cells.Item(31, 1).Value = 20627000
$ws.Cells.Item(31, 2).Value = "Новохопёрский муниципальный район"
$ws.Cells.Item(31, 3).Value = "male"
$ws.Cells.Item(31, 4).Value = 2012
$ws.Cells.Item(31, 5).Value = 0.06128
$ws.Cells.Item(31, 6).Value = 0.05936
$ws.Cells.Item(31, 7).Value = 0.0649
$ws.Cells.Item(31, 8).Value = 0.0636
$ws.Cells.Item(31, 9).Value = 0.0714
$ws.Cells.Item(31, 10).Value = 0.0886
$ws.Cells.Item(31, 11).Value = 0.0767
$ws.Cells.Item(31, 12).Value = 0.0763
$ws.Cells.Item(31, 13).Value = 0.0743
$ws.Cells.Item(31, 14).Value = 0.0863
$ws.Cells.Item(31, 15).Value = 0.10297
$ws.Cells.Item(31, 16).Value = 0.0807
$ws.Cells.Item(31, 17).Value = 0.06665
$ws.Cells.Item(31, 18).Value = 0.02698

# Update the active selection to match the author's final cursor position.
$ws.Range("B33").Select() | Out-Null
